{"js": "// Update existing statistic cells in the \"t1_uninsured\" table, then\n// append the new comorbidity rows at the end of the table.\n\nconst replacements = [\n  [\"36.45 +/- 0.06\", \"34.16 +/- 0.07\"],\n  [\"484 (1.21)\", \"484 (1.10)\"],\n  [\"1,883 (4.70)\", \"1,979 (4.49)\"],\n  [\"664 (1.66)\", \"678 (1.54)\"],\n  [\"15,475 (38.66)\", \"17,001 (38.57)\"],\n  [\"24,385 (60.91)\", \"26,859 (60.93)\"],\n  [\"173 (0.43)\", \"222 (0.50)\"],\n  [\"1,020 (2.55)\", \"1,107 (2.51)\"],\n  [\"4,194 (10.48)\", \"4,391 (9.96)\"],\n  [\"13,494 (33.71)\", \"15,138 (34.34)\"],\n  [\"213 (0.53)\", \"239 (0.54)\"],\n  [\"2,157 (5.39)\", \"2,397 (5.44)\"],\n  [\"1,911 (4.77)\", \"2,240 (5.08)\"],\n  [\"17,044 (42.57)\", \"18,570 (42.13)\"],\n  [\"40,033 (100.00)\", \"44,082 (100.00)\"],\n  [\"4,698 (11.74)\", \"5,108 (11.59)\"],\n  [\"15,200 (37.97)\", \"16,195 (36.74)\"],\n  [\"20,135 (50.30)\", \"22,779 (51.67)\"],\n  [\"5,816 (14.53)\", \"6,413 (14.55)\"],\n  [\"5,398 (13.48)\", \"5,932 (13.46)\"],\n  [\"21,211 (52.98)\", \"23,044 (52.28)\"],\n  [\"7,608 (19.00)\", \"8,693 (19.72)\"],\n  [\"14,830 (37.04)\", \"16,117 (36.56)\"],\n  [\"10,937 (27.32)\", \"12,094 (27.44)\"],\n  [\"8,942 (22.34)\", \"9,908 (22.48)\"],\n  [\"5,324 (13.30)\", \"5,963 (13.53)\"],\n  [\"178 (0.44)\", \"204 (0.46)\"],\n  [\"6,002 (14.99)\", \"6,879 (15.61)\"],\n  [\"9,227 (23.05)\", \"9,960 (22.59)\"],\n];\n\nconst newRows = [\n  [\"AIDS\", \"40.0 (0.09)\"],\n  [\"ALCOHOL\", \"718.0 (1.63)\"],\n  [\"ARTHRITIS\", \"110.0 (0.25)\"],\n  [\"CANCER, LYMPHOMA\", \"11.0 (0.02)\"],\n  [\"CANCER, METASTATIC\", \"26.0 (0.06)\"],\n  [\"CANCER, SOLID\", \"62.0 (0.14)\"],\n  [\"DEPRESSION\", \"896.0 (2.03)\"],\n  [\"DIABETES, UNCOMPLICATED\", \"1,728.0 (3.92)\"],\n  [\"DIABETES, COMPLICATED\", \"212.0 (0.48)\"],\n  [\"HYPERTENSION 1\", \"4,250.0 (9.64)\"],\n  [\"HYPERTENSION 2\", \"3,908.0 (8.87)\"],\n  [\"CHRONIC LUNG\", \"1,850.0 (4.20)\"],\n  [\"OBESITY\", \"2,959.0 (6.71)\"],\n  [\"PERIPHERAL VASCULAR\", \"102.0 (0.23)\"],\n  [\"HYPOTHYROIDISM\", \"546.0 (1.24)\"],\n];\n\nconst body = context.document.body;\n\n// --- Step 1: update the existing value cells in place -----------------\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  found.items[0].insertText(newText, \"Replace\");\n}\nawait context.sync();\n\n// --- Step 2: append the new comorbidity rows to the table -------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.addRows(\"End\", newRows.length, newRows);\nawait context.sync();\n", "ps1": "# Update existing statistic cells in the \"t1_uninsured\" table, then\n# append the new comorbidity rows at the end of the table.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"36.45 +/- 0.06\", \"34.16 +/- 0.07\"),\n    @(\"484 (1.21)\", \"484 (1.10)\"),\n    @(\"1,883 (4.70)\", \"1,979 (4.49)\"),\n    @(\"664 (1.66)\", \"678 (1.54)\"),\n    @(\"15,475 (38.66)\", \"17,001 (38.57)\"),\n    @(\"24,385 (60.91)\", \"26,859 (60.93)\"),\n    @(\"173 (0.43)\", \"222 (0.50)\"),\n    @(\"1,020 (2.55)\", \"1,107 (2.51)\"),\n    @(\"4,194 (10.48)\", \"4,391 (9.96)\"),\n    @(\"13,494 (33.71)\", \"15,138 (34.34)\"),\n    @(\"213 (0.53)\", \"239 (0.54)\"),\n    @(\"2,157 (5.39)\", \"2,397 (5.44)\"),\n    @(\"1,911 (4.77)\", \"2,240 (5.08)\"),\n    @(\"17,044 (42.57)\", \"18,570 (42.13)\"),\n    @(\"40,033 (100.00)\", \"44,082 (100.00)\"),\n    @(\"4,698 (11.74)\", \"5,108 (11.59)\"),\n    @(\"15,200 (37.97)\", \"16,195 (36.74)\"),\n    @(\"20,135 (50.30)\", \"22,779 (51.67)\"),\n    @(\"5,816 (14.53)\", \"6,413 (14.55)\"),\n    @(\"5,398 (13.48)\", \"5,932 (13.46)\"),\n    @(\"21,211 (52.98)\", \"23,044 (52.28)\"),\n    @(\"7,608 (19.00)\", \"8,693 (19.72)\"),\n    @(\"14,830 (37.04)\", \"16,117 (36.56)\"),\n    @(\"10,937 (27.32)\", \"12,094 (27.44)\"),\n    @(\"8,942 (22.34)\", \"9,908 (22.48)\"),\n    @(\"5,324 (13.30)\", \"5,963 (13.53)\"),\n    @(\"178 (0.44)\", \"204 (0.46)\"),\n    @(\"6,002 (14.99)\", \"6,879 (15.61)\"),\n    @(\"9,227 (23.05)\", \"9,960 (22.59)\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n\n$newRows = @(\n    @(\"AIDS\", \"40.0 (0.09)\"),\n    @(\"ALCOHOL\", \"718.0 (1.63)\"),\n    @(\"ARTHRITIS\", \"110.0 (0.25)\"),\n    @(\"CANCER, LYMPHOMA\", \"11.0 (0.02)\"),\n    @(\"CANCER, METASTATIC\", \"26.0 (0.06)\"),\n    @(\"CANCER, SOLID\", \"62.0 (0.14)\"),\n    @(\"DEPRESSION\", \"896.0 (2.03)\"),\n    @(\"DIABETES, UNCOMPLICATED\", \"1,728.0 (3.92)\"),\n    @(\"DIABETES, COMPLICATED\", \"212.0 (0.48)\"),\n    @(\"HYPERTENSION 1\", \"4,250.0 (9.64)\"),\n    @(\"HYPERTENSION 2\", \"3,908.0 (8.87)\"),\n    @(\"CHRONIC LUNG\", \"1,850.0 (4.20)\"),\n    @(\"OBESITY\", \"2,959.0 (6.71)\"),\n    @(\"PERIPHERAL VASCULAR\", \"102.0 (0.23)\"),\n    @(\"HYPOTHYROIDISM\", \"546.0 (1.24)\")\n)\n\n$table = $d.Tables(1)\nforeach ($pair in $newRows) {\n    $row = $table.Rows.Add()\n    $row.Cells(1).Range.Text = $pair[0]\n    $row.Cells(2).Range.Text = $pair[1]\n}\n"}
